# Populate the previously-blank template rows 154-158 of the
# MOSIP_Feature_Roadmap sheet with the new roadmap items described in the
# commit ("Updated Features Roadmap for IDA").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_Feature_Roadmap")

# ---------------------------------------------------------------------
# Row 154 - MOS-21582
# ---------------------------------------------------------------------
$ws.Range("A154").Value = 151
$ws.Range("B154").Value = "MOS-21582"
$ws.Range("C154").NumberFormat = "d-mmm-yy"
$ws.Range("C154").Value = "3/26/2019"
$ws.Range("D154").Value = "API Specification Changes for IDA based on MDS review by Sasi/Ramesh"
$ws.Range("E154").Value = "ID-Authentication"
$ws.Range("F154").Value = "New"
$ws.Range("G154").Value = "Additional or Modification of attributes in API Specs based on review "
$ws.Range("L154").Value = 1
$ws.Range("M154").Value = "Approved"
$ws.Range("N154").Value = "Ramesh"
$ws.Range("O154").NumberFormat = "d-mmm-yy"
$ws.Range("O154").Value = "3/26/2019"

# ---------------------------------------------------------------------
# Row 155 - MOS-21583
# ---------------------------------------------------------------------
$ws.Range("A155").Value = 152
$ws.Range("B155").Value = "MOS-21583"
$ws.Range("C155").NumberFormat = "d-mmm-yy"
$ws.Range("C155").Value = "3/26/2019"
$ws.Range("D155").Value = "Design Change of ID-Repo based on Security review by Sasi/Ramesh"
$ws.Range("E155").Value = "ID-Authentication"
$ws.Range("F155").Value = "New"
$ws.Range("G155").Value = "Design Change of ID-Repo based on Security review by Sasi/Ramesh"
$ws.Range("L155").Value = 1
$ws.Range("M155").Value = "Approved"
$ws.Range("N155").Value = "Ramesh"
$ws.Range("O155").NumberFormat = "d-mmm-yy"
$ws.Range("O155").Value = "3/26/2019"

# ---------------------------------------------------------------------
# Row 156 - MOS-21584
# ---------------------------------------------------------------------
$ws.Range("A156").Value = 153
$ws.Range("B156").Value = "MOS-21584"
$ws.Range("C156").NumberFormat = "d-mmm-yy"
$ws.Range("C156").Value = "3/26/2019"
$ws.Range("D156").Value = "Design Change of IDA based on Security review by Sasi/Ramesh"
$ws.Range("E156").Value = "ID-Authentication"
$ws.Range("F156").Value = "New"
$ws.Range("G156").Value = "Design Change of IDA based on Security review by Sasi/Ramesh"
$ws.Range("L156").Value = 1
$ws.Range("M156").Value = "Approved"
$ws.Range("N156").Value = "Ramesh"
$ws.Range("O156").NumberFormat = "d-mmm-yy"
$ws.Range("O156").Value = "3/26/2019"

# ---------------------------------------------------------------------
# Row 157 - MOS-21585
# ---------------------------------------------------------------------
$ws.Range("A157").Value = 154
$ws.Range("B157").Value = "MOS-21585"
$ws.Range("C157").NumberFormat = "d-mmm-yy"
$ws.Range("C157").Value = "4/1/2019"
$ws.Range("D157").Value = "Mapping of platform address attributes in IDA based on Morrocco Address Structure"
$ws.Range("E157").Value = "ID-Authentication"
$ws.Range("F157").Value = "New"
$ws.Range("G157").Value = "Mapping of platform address attributes in IDA based on Morrocco Address Structure"
$ws.Range("L157").Value = 1
$ws.Range("M157").Value = "Approved"
$ws.Range("N157").Value = "Shrikant"
$ws.Range("O157").NumberFormat = "d-mmm-yy"
$ws.Range("O157").Value = "4/1/2019"

# ---------------------------------------------------------------------
# Row 158 - MOS-21327
# ---------------------------------------------------------------------
$ws.Range("A158").Value = 155
$ws.Range("B158").Value = "MOS-21327"
$ws.Range("C158").NumberFormat = "d-mmm-yy"
$ws.Range("C158").Value = "3/28/2019"
$ws.Range("D158").Value = "Integrate with new VID Generator API"
$ws.Range("E158").Value = "ID-Authentication"
$ws.Range("F158").Value = "New"
$ws.Range("G158").Value = "Integrate with the new VID generator component based on the VID policy/type defined"
$ws.Range("L158").Value = 1
$ws.Range("M158").Value = "Approved"
$ws.Range("N158").Value = "Ramesh"
$ws.Range("O158").NumberFormat = "d-mmm-yy"
$ws.Range("O158").Value = "3/28/2019"

# Move the frozen-pane scroll position / active selection down to the newly
# edited area, mirroring where the author's cursor ended up.
$ws.Range("D156").Select()
